# Re-grade students 48966 and 48982 with plugin doc partial credit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade Comparison")

# Student 48966 -> row 9
$ws.Range("B9").Value = 79
$ws.Range("C9").Value = 64.91220900398571
$ws.Range("D9").Value = -14.08779099601429
$ws.Range("F9").Value = 14.08779099601428

# Student 48982 -> row 18
$ws.Range("B18").Value = 80
$ws.Range("C18").Value = 60.64776531782864
$ws.Range("D18").Value = -19.35223468217136
$ws.Range("F18").Value = 19.35223468217136
$ws.Range("G18").Value = "Good"
